# Add the "Remove Columns" worksheet to the end of the workbook (GOMS sheet
# for the new "Remove Columns" python-IDE task) and populate it with the
# Action / Time / Content walkthrough, mirroring the layout used by the
# other task sheets (e.g. "Datetime Components").

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Remove Columns"

# Header row
$newSheet.Range("A1").Value = "Action"
$newSheet.Range("B1").Value = "Time"
$newSheet.Range("C1").Value = "Content"

# Step rows
$newSheet.Range("A2").Value = "Upload CSV"
$newSheet.Range("B2").Value = "5 min"
$newSheet.Range("C2").Value = "df = pd.read_csv('file.csv')"

$newSheet.Range("A3").Value = "Identify Columns"
$newSheet.Range("B3").Value = "1 min"
$newSheet.Range("C3").Value = "df.columns to list all columns"

$newSheet.Range("A4").Value = "Drop Columns"
$newSheet.Range("B4").Value = "1 min"
$newSheet.Range("C4").Value = "df.drop(['column1', 'column2'], axis=1, inplace=True)"

$newSheet.Range("A5").Value = "Verify Changes"
$newSheet.Range("B5").Value = "1 min"
$newSheet.Range("C5").Value = "df.head() to ensure columns are dropped"

# Overall summary row
$newSheet.Range("A6").Value = "Overall"
$newSheet.Range("B6").Value = "8 min"

# Formatting: header + overall rows bold 13pt, body rows regular 13pt
# (size before bold keeps the engine reusing the workbook's existing
# "bold 13pt" / "regular 13pt" styles instead of minting new ones)
$newSheet.Range("A1:C1").Font.Size = 13
$newSheet.Range("A1:C1").Font.Bold = $true

$newSheet.Range("A2:C5").Font.Size = 13

$newSheet.Range("A6:B6").Font.Size = 13
$newSheet.Range("A6:B6").Font.Bold = $true

# Make the new sheet the active tab/selection, matching the author's
# final view state.
[void]$newSheet.Activate()
[void]$newSheet.Range("G15").Select()
